# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.537.72'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '1.472.46'
$ws.Range('E3').Value = '  +2.28%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.9608'
$ws.Range('E5').Value = '  +4.96%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '277.29'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3598'
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3089'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.083'
$ws.Range('E9').Value = '  +5.78%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.48'
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06637'
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.491'
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.16'
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.167'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9603'
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001022'
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('D18').Value = '1.471.39'
$ws.Range('E18').Value = '  +2.47%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.05947'
$ws.Range('E19').Value = '  +5.82%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '68.99'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.492'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.58'
$ws.Range('E22').Value = '  +2.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.26'
$ws.Range('E23').Value = '  +3.76%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.280'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').Value = '20.565.60'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '144.34'
$ws.Range('E26').Value = '  +4.47%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.115'
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.16'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Value = '1.631.20'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '113.92'
$ws.Range('E30').Value = '  +3.85%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.897'
$ws.Range('E31').Value = '  +1.60%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.8112'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.950'
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07992'
$ws.Range('E34').Value = '  +4.77%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.225'
$ws.Range('E35').Value = '  +8.60%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.466'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05798'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.726'
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02048'
$ws.Range('E39').Value = '  +3.19%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9601'
$ws.Range('E40').Value = '  +4.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '10.41'
$ws.Range('E41').Value = '  +2.49%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1873'
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.448'
$ws.Range('E43').Value = '  +3.74%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5272'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.515'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.15'
$ws.Range('E46').Value = '  +2.39%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '119.06'
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5209'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.811'
$ws.Range('E49').Value = '  +4.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06456'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9880'
$ws.Range('E51').Value = '  +0.09%  '
